# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (F column) numbers across the 展览, 演出 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsPerformance = $wb.Worksheets.Item("演出")
$wsAllTypes = $wb.Worksheets.Item("全部类型")

# ---- 展览 (Exhibition) sheet ----
$wsExhibition.Range("F2").Value = 12904
$wsExhibition.Range("F6").Value = 335
$wsExhibition.Range("F7").Value = 412
$wsExhibition.Range("F9").Value = 13108
$wsExhibition.Range("F10").Value = 49
$wsExhibition.Range("F11").Value = 40
$wsExhibition.Range("F12").Value = 5320
$wsExhibition.Range("F13").Value = 553
$wsExhibition.Range("F14").Value = 24
$wsExhibition.Range("F18").Value = 46
$wsExhibition.Range("F19").Value = 138
$wsExhibition.Range("F20").Value = 694
$wsExhibition.Range("F21").Value = 2868
$wsExhibition.Range("F22").Value = 6243
$wsExhibition.Range("F23").Value = 1169
$wsExhibition.Range("F24").Value = 3648
$wsExhibition.Range("F25").Value = 223
$wsExhibition.Range("F26").Value = 50

# ---- 演出 (Performance) sheet ----
$wsPerformance.Range("F3").Value = 11

# ---- 全部类型 (All Types) sheet ----
$wsAllTypes.Range("F2").Value = 12904
$wsAllTypes.Range("F6").Value = 335
$wsAllTypes.Range("F8").Value = 412
$wsAllTypes.Range("F10").Value = 13108
$wsAllTypes.Range("F11").Value = 49
$wsAllTypes.Range("F12").Value = 40
$wsAllTypes.Range("F13").Value = 5320
$wsAllTypes.Range("F14").Value = 553
$wsAllTypes.Range("F15").Value = 24
$wsAllTypes.Range("F19").Value = 46
$wsAllTypes.Range("F20").Value = 138
$wsAllTypes.Range("F21").Value = 694
$wsAllTypes.Range("F22").Value = 2868
$wsAllTypes.Range("F23").Value = 11
$wsAllTypes.Range("F24").Value = 6243
$wsAllTypes.Range("F25").Value = 1169
$wsAllTypes.Range("F26").Value = 3648
$wsAllTypes.Range("F27").Value = 223
$wsAllTypes.Range("F28").Value = 50

$wb.Save()
